# Applies the cryptos.xlsx price/volume refresh described in the commit:
# "Updated cryptos list on Fri Jul 21 05:22:23 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.863.47'
$ws.Range('E2').Value = '  -0.51%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.895.27'
$ws.Range('E3').Value = '  -0.21%  '

# Row 4
$ws.Range('E4').Value = '  +0.16%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7909'
$ws.Range('E5').Value = '  -5.48%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.61'
$ws.Range('E6').Value = '  +0.63%  '

# Row 7
$ws.Range('E7').Value = '  +0.14%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3153'
$ws.Range('E8').Value = '  -3.95%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.37'
$ws.Range('E9').Value = '  -4.45%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07189'
$ws.Range('E10').Value = '  +2.03%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08105'
$ws.Range('E11').Value = '  +0.25%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.965.30'
$ws.Range('E12').Value = '  +3.48%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.579'
$ws.Range('E13').Value = '  +6.07%  '

# Row 14
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7661'
$ws.Range('E14').Value = '  +0.33%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.52'
$ws.Range('E15').Value = '  +0.15%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.177'
$ws.Range('E16').Value = '  +5.46%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.874.07'
$ws.Range('E17').Value = '  -0.46%  '

# Row 18
$ws.Range('E18').Value = '  -1.28%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.01'
$ws.Range('E19').Value = '  -0.05%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007786'
$ws.Range('E20').Value = '  +0.28%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.198'
$ws.Range('E21').Value = '  +17.70%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.08%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.131.60'
$ws.Range('E23').Value = '  -0.76%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.19%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1657'
$ws.Range('E25').Value = '  -5.06%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.414'
$ws.Range('E26').Value = '  +1.31%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.90'
$ws.Range('E27').Value = '  -0.92%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.70'
$ws.Range('E28').Value = '  -1.26%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.062'
$ws.Range('E29').Value = '  -1.56%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.401'
$ws.Range('E30').Value = '  +2.91%  '

# Row 31
$ws.Range('E31').Value = '  +2.37%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.474'
$ws.Range('E32').Value = '  +4.42%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.091'
$ws.Range('E33').Value = '  +0.44%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05529'
$ws.Range('E34').Value = '  -6.38%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.275'
$ws.Range('E35').Value = '  +0.53%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7399'
$ws.Range('E36').Value = '  +1.02%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9961'
$ws.Range('E37').Value = '  -0.22%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.620'
$ws.Range('E38').Value = '  -3.52%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01927'
$ws.Range('E39').Value = '  +0.49%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.781'
$ws.Range('E40').Value = '  +0.13%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.150.00'
$ws.Range('E41').Value = '  +15.38%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.60'
$ws.Range('E42').Value = '  +2.47%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4423'
$ws.Range('E43').Value = '  -0.53%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.887'
$ws.Range('E44').Value = '  +0.27%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8531'
$ws.Range('E45').Value = '  -0.56%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '104.47'
$ws.Range('E46').Value = '  +2.71%  '

# Row 47
$ws.Range('E47').Value = '  +0.12%  '

# Row 48
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.07'
$ws.Range('E48').Value = '  +2.67%  '

# Row 49
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.877'
$ws.Range('E49').Value = '  -1.60%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.042'
$ws.Range('E50').Value = '  +11.96%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.457'
$ws.Range('E51').Value = '  -1.37%  '
